$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: Restored from revision of admin on 05/31/2021 05:43:42 PM.TEST
# The underlying change is that cell C10's numeric value was changed from 18 to 1.
$ws.Range("C10").Value = 1
